$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @("ac","al","am","ap","ba","ce","df","es","go","ma","mg","ms","mt","pa","pb","pe","pi","pr","rj","rn","ro","rr","rs","sc","se","sp","to")

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

$ws.Range("B1").Select()
